$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price values are plain decimal numbers (e.g. "1.002", "19.93") which Excel would
# otherwise auto-convert to numeric cells. The source data keeps them as text, so we
# pre-format those specific cells as Text before writing the values.
$textCells = @("D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D18", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '27.192.05'
$ws.Range("E2").Value = '  +0.16%  '

# Row 3
$ws.Range("D3").Value = '1.848.82'
$ws.Range("E3").Value = '  +0.11%  '

# Row 4
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.40%  '

# Row 5
$ws.Range("D5").Value = '313.04'
$ws.Range("E5").Value = '  -0.09%  '

# Row 6
$ws.Range("E6").Value = '  -0.34%  '

# Row 7
$ws.Range("D7").Value = '0.4625'
$ws.Range("E7").Value = '  -0.17%  '

# Row 8
$ws.Range("D8").Value = '0.3695'
$ws.Range("E8").Value = '  -0.24%  '

# Row 9
$ws.Range("D9").Value = '0.07264'
$ws.Range("E9").Value = '  -1.48%  '

# Row 10
$ws.Range("D10").Value = '0.8869'
$ws.Range("E10").Value = '  +0.36%  '

# Row 11
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").Value = '19.93'
$ws.Range("E11").Value = '  +0.05%  '

# Row 12
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '0.07833'
$ws.Range("E12").Value = '  -0.97%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.860.86'
$ws.Range("E13").Value = '  +1.00%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '5.391'
$ws.Range("E14").Value = '  +0.34%  '

# Row 15
$ws.Range("D15").Value = '6.503'
$ws.Range("E15").Value = '  -1.19%  '

# Row 16
$ws.Range("D16").Value = '91.38'
$ws.Range("E16").Value = '  -0.55%  '

# Row 17
$ws.Range("E17").Value = '  -0.35%  '

# Row 18
$ws.Range("D18").Value = '0.000008823'
$ws.Range("E18").Value = '  -1.28%  '

# Row 19
$ws.Range("E19").Value = '  -0.31%  '

# Row 20
$ws.Range("D20").Value = '27.229.74'
$ws.Range("E20").Value = '  +0.15%  '

# Row 21
$ws.Range("D21").Value = '14.65'
$ws.Range("E21").Value = '  -1.48%  '

# Row 22
$ws.Range("D22").Value = '5.054'
$ws.Range("E22").Value = '  -1.65%  '

# Row 23
$ws.Range("B23").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D23").Value = '2.113.33'
$ws.Range("E23").Value = '  +0.07%  '

# Row 24
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").Value = '10.53'
$ws.Range("E24").Value = '  -0.38%  '

# Row 25
$ws.Range("D25").Value = '2.037'
$ws.Range("E25").Value = '  +9.31%  '

# Row 26
$ws.Range("D26").Value = '151.46'
$ws.Range("E26").Value = '  -0.97%  '

# Row 27
$ws.Range("D27").Value = '18.38'
$ws.Range("E27").Value = '  -0.56%  '

# Row 28
$ws.Range("D28").Value = '2.022'
$ws.Range("E28").Value = '  -2.21%  '

# Row 29
$ws.Range("D29").Value = '115.53'
$ws.Range("E29").Value = '  -1.28%  '

# Row 30
$ws.Range("D30").Value = '5.003'
$ws.Range("E30").Value = '  -2.41%  '

# Row 31
$ws.Range("D31").Value = '0.08825'
$ws.Range("E31").Value = '  -0.67%  '

# Row 32
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '0.7818'
$ws.Range("E32").Value = '  +5.45%  '

# Row 33
$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").Value = '3.126'
$ws.Range("E33").Value = '  +5.19%  '

# Row 34
$ws.Range("D34").Value = '4.512'
$ws.Range("E34").Value = '  +0.99%  '

# Row 35
$ws.Range("D35").Value = '1.146'
$ws.Range("E35").Value = '  +0.48%  '

# Row 36
$ws.Range("D36").Value = '2.693'
$ws.Range("E36").Value = '  +5.37%  '

# Row 37
$ws.Range("D37").Value = '1.102'
$ws.Range("E37").Value = '  +2.21%  '

# Row 38
$ws.Range("D38").Value = '0.01941'
$ws.Range("E38").Value = '  -0.52%  '

# Row 39
$ws.Range("D39").Value = '0.05206'
$ws.Range("E39").Value = '  -1.20%  '

# Row 40
$ws.Range("D40").Value = '2.949'
$ws.Range("E40").Value = '  -0.70%  '

# Row 41
$ws.Range("D41").Value = '7.015'
$ws.Range("E41").Value = '  -0.91%  '

# Row 42
$ws.Range("D42").Value = '0.5024'
$ws.Range("E42").Value = '  -2.88%  '

# Row 43
$ws.Range("D43").Value = '0.1612'
$ws.Range("E43").Value = '  -1.47%  '

# Row 44
$ws.Range("D44").Value = '8.464'
$ws.Range("E44").Value = '  +2.77%  '

# Row 45
$ws.Range("D45").Value = '0.4751'
$ws.Range("E45").Value = '  -2.22%  '

# Row 46
$ws.Range("D46").Value = '10.34'
$ws.Range("E46").Value = '  +0.69%  '

# Row 47
$ws.Range("E47").Value = '  -0.40%  '

# Row 48
$ws.Range("D48").Value = '103.05'
$ws.Range("E48").Value = '  +0.39%  '

# Row 49
$ws.Range("D49").Value = '1.628'
$ws.Range("E49").Value = '  -0.08%  '

# Row 50
$ws.Range("D50").Value = '0.06201'
$ws.Range("E50").Value = '  -0.50%  '

# Row 51
$ws.Range("D51").Value = '65.51'
$ws.Range("E51").Value = '  -0.19%  '
